# Standardize the "supporting documents" merge-field name: rename every
# occurrence of {supportingDocsList} to {supportingDocs} throughout the
# document body.

$d = $word.ActiveDocument

$wdFindContinue = 1
$wdReplaceAll   = 2

$d.Content.Find.Execute(
    "supportingDocsList",  # FindText
    $true,                 # MatchCase
    $false,                # MatchWholeWord
    $false,                # MatchWildcards
    $false,                # MatchSoundsLike
    $false,                # MatchAllWordForms
    $true,                 # Forward
    $wdFindContinue,       # Wrap
    $false,                # Format
    "supportingDocs",      # ReplaceWith
    $wdReplaceAll          # Replace
) | Out-Null
